# Apply the text-content edits described by the commit's diff.
# (Attribute-order churn on a14:hiddenLine/a14:hiddenFill elements and
#  the presentation-level sldGuideLst ext are serializer artifacts of
#  PowerPoint re-saving the file and are not hand-authored here.)

$p = $ppt.ActivePresentation

# --- Slide 3 ("Key insights") -------------------------------------------
$s3 = $p.Slides.Item(3)
$ph3 = $s3.Shapes.Item(2).TextFrame.TextRange

$ph3.Paragraphs(2).Text = "INC has the unique distinction of contesting every seat in all electoral periods. Similarly, IND has contested every seat during both the Medieval and Latest periods.`r"
$ph3.Paragraphs(3).Text = "TDP maintained a significant presence, contesting the majority of seats in both the medieval and the latest electoral periods. Meanwhile, BJP substantially increased its participation from the medieval to the modern era.`r"

# --- Slide 4 ("Key insights") -------------------------------------------
$s4 = $p.Slides.Item(4)
$ph4 = $s4.Shapes.Item(2).TextFrame.TextRange

$ph4.Paragraphs(2).Text = "Over time, the INC experienced a decline in its vote share, yet it continued to secure the highest percentage of votes. Conversely, IND saw a consistent decrease in their vote share, while TDP witnessed an upward trend in theirs.`r"

# --- Slide 5 ("Final conclusion") ---------------------------------------
$s5 = $p.Slides.Item(5)
$sh5 = $s5.Shapes.Item(2)
$ph5 = $sh5.TextFrame.TextRange

$ph5.Paragraphs(1).Text = " Adult candidates have gathered a larger proportion of votes relative to their representation among candidates. Moreover, fielding candidates in every seat does not necessarily translate to a higher vote count; for instance, IND contested the majority of seats over time, yet their share of votes diminished. Conversely, TDP ranked fourth in seat participation but secured the second-highest vote share during the medieval period. In the most recent period, TDP's participation was the third-highest for GEN seats and second-highest for SC & ST seats, yet it outperformed the INC in vote share, claiming the top spot.`r"

$sh5.TextFrame.AutoSize = $ppAutoSizeNone
$sh5.TextFrame2.TextRange.Font.Fill.ForeColor.RGB = $sh5.TextFrame2.TextRange.Font.Fill.ForeColor.RGB
